# Week 5 final commit
# Applies the schedule content changes: edits existing rows 3-17 (text
# tweaks / corrected labels), fills in the brand-new Week 5 ("Brussels
# rent prices") row that got inserted, and populates weeks 17-21 which
# were previously blank placeholders.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- Row 4 (Week 3): Goal / Keywords swapped in for the crack-detection task ---
$ws.Range("D4").Value = "Compare how fast and how good all methods are, including how easy to use"
$ws.Range("E4").Value = "Supervised Learning, Gradient Boosting, XGBoost, LightGBM, CatBoost, K-Fold Cross Validation, Hyperparameter Tuning"

# --- Row 5 (Week 4): Data Source / Goal text refreshed ---
$ws.Range("C5").Value = "https://www.kaggle.com/datasets/sinamhd9/concrete-comprehensive-strength?select=Concrete_Data.xls"
$ws.Range("D5").Value = "Build a predictive model"

# --- Row 6 (Week 5): brand-new "Brussels rent prices" entry ---
$ws.Range("B6").Value = "Brussels Rent Prices Webscraping Selenium + Modelling"
$ws.Range("C6").Value = "https://www.immoweb.be/en/search/house-and-apartment/for-rent/brussels/province?countries=BE&page=1&orderBy=newest"
$ws.Range("D6").Value = "Scrape rent prices for all of Brussels from the website of a major Belgian real estate platform. Then clean, visualise and model these data."
$ws.Range("E6").Value = "Web Scraping, Selenium, RSelenium, ggplot2, dplyr, pandas, seaborn, matplotlib, Feature Selection, Data Cleaning, Regression Models"

# --- Row 7 (Week 6): title clarified ---
$ws.Range("B7").Value = "Scraping and visualising housing prices for different post codes in Vienna, including descriptions"

# --- Row 8 (Week 7): title reworded ---
$ws.Range("B8").Value = "Predicting used car prices"

# --- Rows 9-17 (Weeks 8-16): same titles, shifted down from their old slots ---
$ws.Range("B9").Value = "Loans (Data is Plural)"
$ws.Range("B10").Value = "Procrastinated one about importance of retraining model (Day ahead spot prices)"
$ws.Range("B11").Value = "Image classification: Cracks in concrete"
$ws.Range("C11").Value = "https://www.kaggle.com/code/vishnu0399/ensuring-structural-safety-crack-detection"
$ws.Range("B12").Value = "Electricity something"
$ws.Range("B13").Value = "Prisoners problem"
$ws.Range("B14").Value = "Pytorch"
$ws.Range("B15").Value = "Webscraping Digitec or Ebay Tech Prices (Focus on one, either smartphones, laptops or something else, maybe Apple)"
$ws.Range("B16").Value = "F1 2022 vs 2021 data?"
$ws.Range("B17").Value = "Luxembourg Stats (Compare weather with St. Gallen?)"

# --- Rows 18-22 (Weeks 17-21): previously empty, now populated ---
$ws.Range("B18").Value = "Replicate David Robinson Video"
$ws.Range("B19").Value = "World Bank Data?"
$ws.Range("B20").Value = "Bern Rent Prices vs. Vienna vs. Brussels vs. London vs. New York"
$ws.Range("B21").Value = "New York Open Data"
$ws.Range("C21").Value = "https://data.cityofnewyork.us/browse?category=Health"
$ws.Range("B22").Value = "Scraping Rolex Prices"

# --- Row 49 (Week 48) / Row 52 (Week 51): unchanged text, just re-pointed ---
$ws.Range("B49").Value = "Unisport Bern Occupancy"
$ws.Range("B52").Value = "Lookback on the Gym Year (Gym Tracking)"

# --- View state: scroll down a bit and select C22, matching the author's
#     last on-screen position when they saved. ---
$ws.Range("C22").Select()
$win = $excel.ActiveWindow
$win.ScrollRow = 10
$win.ScrollColumn = 1
